# Commit: swap the "Office Theme" and "Integral" theme definitions so the
# deck's live theme (ppt/theme/theme2.xml, wired to the slide master /
# notes master / presentation) carries the stock "Office Theme" palette
# that used to live in ppt/theme/theme1.xml (and vice versa conceptually).
#
# The PowerPoint object model only lets automation rewrite the *colour
# values* of the active theme's colour scheme (Design.SlideMaster's
# ColorScheme / a Slide's ThemeColorScheme) - the theme/clrScheme display
# names and the unused theme1.xml part are not reachable through COM, so
# we reproduce the reachable, meaningful part of the change: restoring
# the master theme's 12 scheme colours to the Office Theme defaults.

$p = $ppt.ActivePresentation
$master = $p.Designs.Item(1).SlideMaster
$cs = $master.ColorScheme

function HexToOleRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme scheme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# - these used to live in ppt/theme/theme1.xml and now become the colours
# used by the live/master theme (ppt/theme/theme2.xml).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000",  # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $cs.Colors($i + 1).RGB = HexToOleRgb($officeThemeColors[$i])
}
